$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1 and copy the formatting used by the other
# header cells (e.g. G1) so it matches style index 1 (bold, bordered,
# centered header style).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Add the corresponding data value for the new "Save" column in row 2
$ws.Range("H2").Value = 0
